# Convert Phone (D) and Mobile (E) columns from text-numbers to true numbers
# for rows 2-91, fix row 92 (dedupe cleanup), and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$phoneMobileNumbers = @{
    "D2" = 2731026620
    "E2" = 6974754853
    "D3" = 2155309555
    "E3" = 6947695629
    "D4" = 2130333815
    "E4" = 6944744724
    "D5" = 2106036036
    "E5" = 6974121315
    "D6" = 2105902059
    "D7" = 2130316283
    "E7" = 6943828282
    "D8" = 2109840078
    "E8" = 6977202619
    "D9" = 2105767889
    "E9" = 6973492279
    "D10" = 2106982208
    "E10" = 6944911408
    "D11" = 2105025268
    "E11" = 6972070111
    "D12" = 2106148852
    "E12" = 6987578684
    "D13" = 2114013841
    "E13" = 6932441469
    "D14" = 2105713710
    "E14" = 6945404806
    "D15" = 2107222551
    "E15" = 6972503355
    "D16" = 2114180300
    "E16" = 6932225937
    "D17" = 2105905460
    "E17" = 6936988762
    "D18" = 2102406299
    "E18" = 6942076290
    "D19" = 2106922953
    "E19" = 6945647623
    "D20" = 2109248838
    "E20" = 6972202914
    "D21" = 2103455493
    "E21" = 6932351230
    "D22" = 2109645820
    "E22" = 6942820500
    "D23" = 2104111821
    "E23" = 6944722872
    "D24" = 2294023771
    "E24" = 6936608099
    "D25" = 2103622319
    "E25" = 6944962722
    "D26" = 2102917955
    "E26" = 6974533639
    "D27" = 2108079921
    "E27" = 6944698898
    "D28" = 2106753795
    "E28" = 6974484918
    "D29" = 2106745205
    "E29" = 6980351683
    "D30" = 2810220570
    "E30" = 6947328410
    "D31" = 2106890513
    "D32" = 2106834141
    "D33" = 2106000611
    "D34" = 2108002495
    "E34" = 6977325007
    "D35" = 2106856106
    "D36" = 2107717107
    "D37" = 2106746264
    "E37" = 6941499153
    "D38" = 2106838647
    "E38" = 6972246373
    "D39" = 2106756420
    "D40" = 2108311063
    "E40" = 6974115341
    "D41" = 2106754422
    "D42" = 2106000757
    "D43" = 2106826514
    "E43" = 6944281011
    "D44" = 2106724944
    "D45" = 2102838196
    "E45" = 6978875087
    "D46" = 2102843353
    "D47" = 2106817624
    "D48" = 2106835741
    "E48" = 6945592521
    "D49" = 2106515029
    "E49" = 6948087946
    "D50" = 2109566774
    "E50" = 6970506798
    "D51" = 2111158777
    "E51" = 6978222644
    "D52" = 2111181200
    "E52" = 6981033171
    "D53" = 2102449412
    "E53" = 6974953295
    "D54" = 2109315648
    "E54" = 6974345938
    "D55" = 2109591846
    "E55" = 6945857491
    "D56" = 2295054665
    "D57" = 2114103647
    "E57" = 6972838250
    "D58" = 2109919010
    "E58" = 6944414808
    "D59" = 2121053500
    "E59" = 6937374680
    "D60" = 2105014500
    "E60" = 6977037641
    "D61" = 2109658906
    "E61" = 6932632898
    "D62" = 2108225065
    "E62" = 6947434884
    "D63" = 2102799540
    "E63" = 6946369199
    "D64" = 2107799566
    "E64" = 6944654891
    "D65" = 2114081622
    "E65" = 6947326295
    "D66" = 2168081782
    "E66" = 6945345458
    "D67" = 2106610209
    "E67" = 6974016795
    "D68" = 2111165089
    "E68" = 6941505426
    "D69" = 2299049199
    "E69" = 6977695222
    "D70" = 2102855326
    "E70" = 6977562551
    "D71" = 6944542221
    "D72" = 2108829023
    "E72" = 6974661171
    "D73" = 2106124549
    "E73" = 6942409460
    "D74" = 2109646114
    "D75" = 2102930460
    "E75" = 6944281772
    "D76" = 2109646400
    "E76" = 6944348562
    "D77" = 2102014314
    "E77" = 6972190704
    "D78" = 2108830727
    "E78" = 6979494266
    "D79" = 2107482041
    "E79" = 6972190704
    "D80" = 2106148114
    "E80" = 6974492510
    "D81" = 2155255842
    "E81" = 6973531577
    "D82" = 2117350265
    "D83" = 2299063249
    "E83" = 6999440441
    "D84" = 2105057277
    "E84" = 6972076641
    "D85" = 2109833368
    "E85" = 6978061340
    "D86" = 2130231610
    "E86" = 6977658236
    "D87" = 2104940317
    "E87" = 6978613054
    "D88" = 2109414967
    "E88" = 6944900690
    "D89" = 2107665604
    "E89" = 6946246807
    "D90" = 2107774666
    "E90" = 6972920788
    "D91" = 6944197808
}

foreach ($cellRef in $phoneMobileNumbers.Keys) {
    $ws.Range($cellRef).Value = $phoneMobileNumbers[$cellRef]
}

# Row 92 cleanup: replace the doctor name with "egertg", move the remaining
# mobile number into D92 as a real number, and clear the old E92 value.
$ws.Range("A92").Value = "egertg"
$ws.Range("D92").Value = 6944197808
$ws.Range("E92").ClearContents()

# Restore the active selection/scroll position reported in the saved view
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 82
$win.ScrollColumn = 1
$ws.Range("E93").Select() | Out-Null

Write-Host "Edit complete"
